$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Matteo Mazzola"
$ws.Range("B5").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C5").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("D5").Value = "Daniele Feller | GREP"
$ws.Range("E5").Value = "Federico  Zanini | A.C.DENTI"
$ws.Range("F5").Value = "Luca Giordani | SHARK ATTACK"
